$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5010.357
$ws.Range("I62").Value = 2214.5
$ws.Range("J62").Value = 12000
$ws.Range("K62").Value = 2214.5
$ws.Range("L62").Value = 12000
$ws.Range("M62").Value = -1590.5
$ws.Range("N62").Value = -13248
$ws.Range("H65").Value = 5010.357
$ws.Range("I65").Value = 2214.5
$ws.Range("J65").Value = 12000
$ws.Range("K65").Value = 11072.5
$ws.Range("L65").Value = 60000
$ws.Range("M65").Value = -7952.5
$ws.Range("N65").Value = -66240
$ws.Range("H70").Value = 9999.666999999999
$ws.Range("I70").Value = 9999
$ws.Range("K70").Value = 29997
$ws.Range("M70").Value = -29727
$ws.Range("H73").Value = 9999.666999999999
$ws.Range("I73").Value = 9999
$ws.Range("K73").Value = 29997
$ws.Range("M73").Value = -29061
$ws.Range("H116").Value = 5546.6665
$ws.Range("I116").Value = 10787.5
$ws.Range("K116").Value = 10787.5
$ws.Range("M116").Value = -7345.5
$ws.Range("H132").Value = 5817.25
$ws.Range("J132").Value = 24666.666
$ws.Range("L132").Value = 73999.99800000001
$ws.Range("N132").Value = -79059.99800000001

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 11858.842
$ws.Range("J44").Value = 11858.842
$ws.Range("L44").Value = 11858.842
$ws.Range("N44").Value = -12834.842
$ws.Range("H61").Value = 1241.9474
$ws.Range("I61").Value = 1269.3889
$ws.Range("K61").Value = 1269.3889
$ws.Range("M61").Value = -1057.3889
$ws.Range("H74").Value = 4675.95
$ws.Range("I74").Value = 4554.353
$ws.Range("K74").Value = 4554.353
$ws.Range("M74").Value = -3680.353
$ws.Range("H77").Value = 4675.95
$ws.Range("I77").Value = 4554.353
$ws.Range("K77").Value = 22771.765
$ws.Range("M77").Value = -18403.765
$ws.Range("H136").Value = 1241.9474
$ws.Range("I136").Value = 1269.3889
$ws.Range("K136").Value = 3808.1667
$ws.Range("M136").Value = -1258.1667

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 1787.5
$ws.Range("I37").Value = 1395.8572
$ws.Range("K37").Value = 1395.8572
$ws.Range("M37").Value = -1258.8572
$ws.Range("H134").Value = 1550.3334
$ws.Range("I134").Value = 1342.9
$ws.Range("K134").Value = 4028.7
$ws.Range("M134").Value = -1493.7

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1678.3
$ws.Range("I16").Value = 1160.375
$ws.Range("K16").Value = 1160.375
$ws.Range("M16").Value = -873.375
$ws.Range("H70").Value = 70750
$ws.Range("J70").Value = 70750
$ws.Range("L70").Value = 70750
$ws.Range("N70").Value = -71380
$ws.Range("H73").Value = 70750
$ws.Range("J73").Value = 70750
$ws.Range("L73").Value = 70750
$ws.Range("N73").Value = -72934
$ws.Range("H107").Value = 998.7895
$ws.Range("I107").Value = 803
$ws.Range("J107").Value = 1334.4286
$ws.Range("K107").Value = 803
$ws.Range("L107").Value = 1334.4286
$ws.Range("M107").Value = 1117
$ws.Range("N107").Value = -5174.4286
$ws.Range("H113").Value = 1678.3
$ws.Range("I113").Value = 1160.375
$ws.Range("K113").Value = 1160.375
$ws.Range("M113").Value = 1009.625
$ws.Range("H132").Value = 3906.842
$ws.Range("I132").Value = 3651.2144
$ws.Range("J132").Value = 4622.6
$ws.Range("K132").Value = 10953.6432
$ws.Range("L132").Value = 13867.8
$ws.Range("M132").Value = -8423.643199999999
$ws.Range("N132").Value = -18927.8

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2526.1428
$ws.Range("J55").Value = 3545.4285
$ws.Range("L55").Value = 10636.2855
$ws.Range("N55").Value = -10990.2855
$ws.Range("H122").Value = 870.1111
$ws.Range("I122").Value = 866.5
$ws.Range("J122").Value = 877.3333
$ws.Range("K122").Value = 7798.5
$ws.Range("L122").Value = 7895.9997
$ws.Range("M122").Value = -5348.5
$ws.Range("N122").Value = -12795.9997

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 34000
$ws.Range("J39").Value = 34000
$ws.Range("L39").Value = 34000
$ws.Range("N39").Value = -35064
$ws.Range("H102").Value = 1627.5454
$ws.Range("I102").Value = 1627.5454
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1627.5454
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -5.545399999999972
$ws.Range("N102").ClearContents()
$ws.Range("H126").Value = 3205.4375
$ws.Range("I126").Value = 2949.0715
$ws.Range("K126").Value = 8847.2145
$ws.Range("M126").Value = -6377.2145
$ws.Range("H132").Value = 53839.45
$ws.Range("I132").Value = 79364.69500000001
$ws.Range("J132").Value = 6435.4287
$ws.Range("K132").Value = 238094.085
$ws.Range("L132").Value = 19306.2861
$ws.Range("M132").Value = -235564.085
$ws.Range("N132").Value = -24366.2861
$ws.Range("H134").Value = 112610.25
$ws.Range("J134").Value = 112610.25
$ws.Range("L134").Value = 337830.75
$ws.Range("N134").Value = -342900.75

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 986.4286
$ws.Range("I31").Value = 533.3333
$ws.Range("J31").Value = 1326.25
$ws.Range("K31").Value = 533.3333
$ws.Range("L31").Value = 1326.25
$ws.Range("M31").Value = -285.3333
$ws.Range("N31").Value = -1822.25
$ws.Range("H32").Value = 9000
$ws.Range("I32").Value = 2500
$ws.Range("K32").Value = 2500
$ws.Range("M32").Value = -2183
$ws.Range("H68").Value = 6278.9
$ws.Range("I68").Value = 4997.3335
$ws.Range("J68").Value = 6828.143
$ws.Range("K68").Value = 4997.3335
$ws.Range("L68").Value = 6828.143
$ws.Range("M68").Value = -4248.3335
$ws.Range("N68").Value = -8326.143
$ws.Range("H71").Value = 6278.9
$ws.Range("I71").Value = 4997.3335
$ws.Range("J71").Value = 6828.143
$ws.Range("K71").Value = 24986.6675
$ws.Range("L71").Value = 34140.715
$ws.Range("M71").Value = -21242.6675
$ws.Range("N71").Value = -41628.715
$ws.Range("H122").Value = 4192.2856
$ws.Range("I122").Value = 4079.4
$ws.Range("K122").Value = 12238.2
$ws.Range("M122").Value = -9788.200000000001
$ws.Range("H136").Value = 2844.4443
$ws.Range("I136").Value = 2385.7144
$ws.Range("K136").Value = 7157.1432
$ws.Range("M136").Value = -4607.1432

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 7500
$ws.Range("I20").Value = 7500
$ws.Range("K20").Value = 7500
$ws.Range("M20").Value = -7260
$ws.Range("H51").Value = 17425
$ws.Range("I51").Value = 17425
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 17425
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -16915
$ws.Range("N51").ClearContents()
$ws.Range("H62").Value = 10857
$ws.Range("I62").Value = 7000
$ws.Range("J62").Value = 11499.833
$ws.Range("K62").Value = 7000
$ws.Range("L62").Value = 11499.833
$ws.Range("M62").Value = -6376
$ws.Range("N62").Value = -12747.833
$ws.Range("H65").Value = 10857
$ws.Range("I65").Value = 7000
$ws.Range("J65").Value = 11499.833
$ws.Range("K65").Value = 35000
$ws.Range("L65").Value = 57499.165
$ws.Range("M65").Value = -31880
$ws.Range("N65").Value = -63739.165
$ws.Range("H75").Value = 40000
$ws.Range("I75").Value = 40000
$ws.Range("K75").Value = 40000
$ws.Range("M75").Value = -39064
$ws.Range("H78").Value = 40000
$ws.Range("I78").Value = 40000
$ws.Range("K78").Value = 120000
$ws.Range("M78").Value = -115320
$ws.Range("H113").Value = 1733
$ws.Range("I113").Value = 1849.6666
$ws.Range("J113").Value = 1499.6666
$ws.Range("K113").Value = 5548.9998
$ws.Range("L113").Value = 4498.9998
$ws.Range("M113").Value = -3378.9998
$ws.Range("N113").Value = -8838.9998
$ws.Range("H122").Value = 3727.2
$ws.Range("I122").Value = 1570.3334
$ws.Range("K122").Value = 4711.0002
$ws.Range("M122").Value = -2261.0002
$ws.Range("H132").Value = 2062.6365
$ws.Range("I132").Value = 2148.9
$ws.Range("J132").Value = 1200
$ws.Range("K132").Value = 6446.700000000001
$ws.Range("L132").Value = 3600
$ws.Range("M132").Value = -3916.700000000001
$ws.Range("N132").Value = -8660
$ws.Range("H136").Value = 4370.273
$ws.Range("I136").Value = 3785.889
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 11357.667
$ws.Range("L136").Value = 21000
$ws.Range("M136").Value = -8807.667000000001
$ws.Range("N136").Value = -26100
